$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''34.127.52'
$ws.Range("E2").Value = '  +0.67%  '
$ws.Range("D3").Value = '''1.783.86'
$ws.Range("E3").Value = '  -1.45%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '''226.42'
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("E7").Value = '  +0.25%  '
$ws.Range("D8").Value = '''31.46'
$ws.Range("E8").Value = '  +1.58%  '
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("D10").Value = '''0.0660'
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = '''2.044.94'
$ws.Range("E12").Value = '  -1.23%  '
$ws.Range("D13").Value = '''11.26'
$ws.Range("E13").Value = '  +11.92%  '
$ws.Range("D14").Value = '''1.792.78'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").Value = '''0.629'
$ws.Range("E15").Value = '  -1.00%  '
$ws.Range("D16").Value = '''34.138.20'
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '''4.22'
$ws.Range("E17").Value = '  -0.32%  '
$ws.Range("D18").Value = '''69.19'
$ws.Range("E18").Value = '  -0.01%  '
$ws.Range("D19").Value = '''254.05'
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").Value = '''0.0₃0741'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '''10.47'
$ws.Range("E22").Value = '  +0.67%  '
$ws.Range("D23").Value = '''4.21'
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("E24").Value = '  -1.66%  '
$ws.Range("D25").Value = '''156.14'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("D26").Value = '''16.48'
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("D27").Value = '''7.04'
$ws.Range("E27").Value = '  +0.24%  '
$ws.Range("E28").Value = '  -0.48%  '
$ws.Range("E29").Value = '  +0.11%  '
$ws.Range("D30").Value = '''3.78'
$ws.Range("E30").Value = '  -0.67%  '
$ws.Range("D31").Value = '''0.0517'
$ws.Range("E31").Value = '  +1.93%  '
$ws.Range("E32").Value = '  +0.29%  '
$ws.Range("D33").Value = '''3.57'
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("E34").Value = '  +2.58%  '
$ws.Range("D35").Value = '''1.451.20'
$ws.Range("E35").Value = '  -5.59%  '
$ws.Range("D36").Value = '''1.06'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = '''0.633'
$ws.Range("E37").Value = '  +2.82%  '
$ws.Range("E38").Value = '  +1.24%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '''2.86'
$ws.Range("E39").Value = '  +1.46%  '
$ws.Range("B40").Value = 'Aave'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D40").Value = '''83.21'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("E41").Value = '  +0.63%  '
$ws.Range("D42").Value = '''0.896'
$ws.Range("E42").Value = '  -0.58%  '
$ws.Range("E43").Value = '  -1.62%  '
$ws.Range("D44").Value = '''0.0510'
$ws.Range("E44").Value = '  -2.05%  '
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '''1.943.23'
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''5.82'
$ws.Range("E47").Value = '  +3.41%  '
$ws.Range("E48").Value = '  +0.22%  '
$ws.Range("D49").Value = '''11.90'
$ws.Range("E49").Value = '  +6.62%  '
$ws.Range("D50").Value = '''50.35'
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("D51").Value = '''97.90'
$ws.Range("E51").Value = '  +2.21%  '
